# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sun Nov 10 10:31:22 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: a literal leading apostrophe forces Excel to keep numeric-looking
# text (e.g. "204.85") as text instead of silently converting it to a number.
$q = "'"

$ws.Range('D2').Value = '79.306.83'
$ws.Range('E2').Value = '  +3.57%  '
$ws.Range('D3').Value = '3.189.06'
$ws.Range('E3').Value = '  +4.62%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = $q + '204.85'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = $q + '633.93'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = $q + '0.234'
$ws.Range('E8').Value = '  +12.33%  '
$ws.Range('D9').Value = $q + '0.584'
$ws.Range('E9').Value = '  +5.80%  '
$ws.Range('D10').Value = '3.188.77'
$ws.Range('E10').Value = '  +4.60%  '
$ws.Range('E11').Value = '  +33.26%  '
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').Value = $q + '5.51'
$ws.Range('E13').Value = '  +7.52%  '
$ws.Range('D14').Value = '3.778.61'
$ws.Range('E14').Value = '  +4.59%  '
$ws.Range('D15').Value = $q + '0.0000227'
$ws.Range('E15').Value = '  +17.49%  '
$ws.Range('D16').Value = $q + '31.82'
$ws.Range('E16').Value = '  +7.20%  '
$ws.Range('D17').Value = '79.223.64'
$ws.Range('E17').Value = '  +3.51%  '
$ws.Range('D18').Value = '3.194.16'
$ws.Range('E18').Value = '  +4.73%  '
$ws.Range('D19').Value = $q + '14.49'
$ws.Range('E19').Value = '  +7.37%  '
$ws.Range('D20').Value = $q + '3.02'
$ws.Range('E20').Value = '  +31.66%  '
$ws.Range('D21').Value = $q + '9.22'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').Value = $q + '426.48'
$ws.Range('E22').Value = '  +13.54%  '
$ws.Range('D23').Value = $q + '5.04'
$ws.Range('E23').Value = '  +15.80%  '
$ws.Range('B24').Value = 'Aptos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D24').Value = $q + '11.37'
$ws.Range('E24').Value = '  +14.23%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.362.12'
$ws.Range('E25').Value = '  +4.84%  '
$ws.Range('D26').Value = $q + '4.77'
$ws.Range('E26').Value = '  +8.56%  '
$ws.Range('D27').Value = $q + '76.95'
$ws.Range('E27').Value = '  +4.54%  '
$ws.Range('D28').Value = $q + '0.998'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = $q + '0.0000118'
$ws.Range('E29').Value = '  +6.01%  '
$ws.Range('D30').Value = $q + '0.997'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  +8.20%  '
$ws.Range('E32').Value = '  +4.06%  '
$ws.Range('D33').Value = $q + '521.65'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('D34').Value = $q + '1.99'
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('E35').Value = '  +28.54%  '
$ws.Range('D36').Value = $q + '22.84'
$ws.Range('E36').Value = '  +9.02%  '
$ws.Range('E37').Value = '  +11.75%  '
$ws.Range('D38').Value = $q + '0.999'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = $q + '0.404'
$ws.Range('E39').Value = '  +4.58%  '
$ws.Range('D40').Value = $q + '165.31'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').Value = $q + '20.03'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = $q + '1.00'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = $q + '191.92'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('D44').Value = $q + '5.51'
$ws.Range('E44').Value = '  +5.60%  '
$ws.Range('D45').Value = $q + '0.813'
$ws.Range('E45').Value = '  +10.00%  '
$ws.Range('E46').Value = '  +7.18%  '
$ws.Range('E47').Value = '  +3.83%  '
$ws.Range('D48').Value = $q + '43.17'
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('D49').Value = $q + '25.71'
$ws.Range('E49').Value = '  +14.55%  '
$ws.Range('D50').Value = $q + '0.635'
$ws.Range('E50').Value = '  +4.49%  '
$ws.Range('D51').Value = $q + '2.50'
$ws.Range('E51').Value = '  +2.40%  '
